$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Simplify the "init" state-table formulas to "none -> none" / "none"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("{ i:Int, p : PID | p != p }", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "none -> none", 2) | Out-Null

$d.Content.Find.Execute("{ p :PID | p ! =p }", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "none", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Remove the stray "_GoBack" bookmark that currently sits at the end
#    of the removePhoto table (it will be re-created further down, at
#    the new end-of-document edit location). Rebuilding the row is the
#    only reliable way to drop the bookmark in this environment.
# ---------------------------------------------------------------------
$removePhotoTable = $d.Tables.Item(6)
$lastRowOld = $removePhotoTable.Rows.Item($removePhotoTable.Rows.Count)
$lastRowOld.Delete()
$removePhotoTable.Rows.Add() | Out-Null
$lastRowIdx = $removePhotoTable.Rows.Count
$removePhotoTable.Cell($lastRowIdx, 1).Range.Text = "toDelete"
$removePhotoTable.Cell($lastRowIdx, 2).Range.Text = "toDelete"
$removePhotoTable.Cell($lastRowIdx, 3).Range.Text = "toDelete + album[i]"

# ---------------------------------------------------------------------
# 3) Declare the new output variable "output_r:Int" on the addPhoto
#    fixture table.
# ---------------------------------------------------------------------
$addPhotoTable = $d.Tables.Item(5)
$pidCell = $addPhotoTable.Cell(3, 1)
$pidCell.Range.Paragraphs.Item(1).Range.InsertAfter(", output_r:Int")

# ---------------------------------------------------------------------
# 4) Add the new "output_r" row (initial value 1, and 0 for every
#    operation) at the bottom of the addPhoto fixture table.
# ---------------------------------------------------------------------
$addPhotoTable.Rows.Add() | Out-Null
$newRowIdx = $addPhotoTable.Rows.Count
$addPhotoTable.Cell($newRowIdx, 1).Range.Text = "output_r"
$addPhotoTable.Cell($newRowIdx, 2).Range.Text = "1"
$addPhotoTable.Cell($newRowIdx, 3).Range.Text = "0"
$addPhotoTable.Cell($newRowIdx, 4).Range.Text = "0"
$addPhotoTable.Cell($newRowIdx, 5).Range.Text = "0"

# Re-create the "_GoBack" bookmark at this new last-edited location.
$lastCell = $addPhotoTable.Cell($newRowIdx, 5)
$lastCellRange = $lastCell.Range.Duplicate
$bookmarkRange = $d.Range($lastCellRange.Start, $lastCellRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
